$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "NAME" column header in B5 (the decision table header row).
# It uses a plain (non-bold) black font with no fill, matching the rest of
# the default table text but with an explicit black font color.
$nameCell = $ws.Range("B5")
$nameCell.Value = "NAME"
$nameCell.Font.Bold = $false
$nameCell.Font.Color = 0

# The rule-name label that used to live in B8 is removed (the column now
# gets its label from the new B5 header instead), leaving B8 blank but
# keeping its existing formatting.
$ws.Range("B8").ClearContents()

# Move the active selection, matching the saved workbook state.
$ws.Range("C16").Select()
